$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1339.3636
$ws.Range("J32").Value = 1365.8889
$ws.Range("L32").Value = 1365.8889
$ws.Range("N32").Value = -2017.8889

$ws.Range("H107").Value = 372.33334
$ws.Range("I107").Value = 205.61111
$ws.Range("K107").Value = 205.61111
$ws.Range("M107").Value = 1714.38889

$ws.Range("H112").Value = 2057.742
$ws.Range("I112").Value = 737.7778
$ws.Range("J112").Value = 2597.7273
$ws.Range("K112").Value = 2213.3334
$ws.Range("L112").Value = 7793.1819
$ws.Range("M112").Value = -1105.3334
$ws.Range("N112").Value = -10009.1819

$ws.Range("H137").Value = 2624.9
$ws.Range("I137").Value = 2694.3333
$ws.Range("K137").Value = 8082.999899999999
$ws.Range("M137").Value = -5532.999899999999

$ws.Range("H138").Value = 2179.2727
$ws.Range("I138").Value = 1820.2354
$ws.Range("J138").Value = 3400
$ws.Range("K138").Value = 5460.706200000001
$ws.Range("L138").Value = 10200
$ws.Range("M138").Value = -320.7062000000005
$ws.Range("N138").Value = -20480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -504
$ws.Range("N97").ClearContents()

$ws.Range("H132").Value = 3891.5862
$ws.Range("I132").Value = 1959.697
$ws.Range("J132").Value = 6441.68
$ws.Range("K132").Value = 5879.090999999999
$ws.Range("L132").Value = 19325.04
$ws.Range("M132").Value = -3349.090999999999
$ws.Range("N132").Value = -24385.04

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2919.1738
$ws.Range("I99").Value = 2583.8462
$ws.Range("J99").Value = 3355.1
$ws.Range("K99").Value = 2583.8462
$ws.Range("L99").Value = 3355.1
$ws.Range("M99").Value = -1085.8462
$ws.Range("N99").Value = -6351.1

$ws.Range("H105").Value = 4508.5713
$ws.Range("I105").Value = 4593.3335
$ws.Range("K105").Value = 4593.3335
$ws.Range("M105").Value = -2846.3335

$ws.Range("H107").Value = 911
$ws.Range("I107").Value = 911
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 911
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1009
$ws.Range("N107").ClearContents()

$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954

$ws.Range("H134").Value = 4532.1963
$ws.Range("I134").Value = 2559
$ws.Range("J134").Value = 6584.32
$ws.Range("K134").Value = 7677
$ws.Range("L134").Value = 19752.96
$ws.Range("M134").Value = -5142
$ws.Range("N134").Value = -24822.96

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6805074.5
$ws.Range("I31").Value = 1687.2
$ws.Range("J31").Value = 17547266
$ws.Range("K31").Value = 1687.2
$ws.Range("L31").Value = 17547266
$ws.Range("M31").Value = -1392.2
$ws.Range("N31").Value = -17547856

$ws.Range("H34").Value = 6805074.5
$ws.Range("I34").Value = 1687.2
$ws.Range("J34").Value = 17547266
$ws.Range("K34").Value = 1687.2
$ws.Range("L34").Value = 17547266
$ws.Range("M34").Value = -1485.2
$ws.Range("N34").Value = -17547670

$ws.Range("H86").Value = 58831052
$ws.Range("I86").Value = 142866990
$ws.Range("J86").Value = 5899.9
$ws.Range("K86").Value = 142866990
$ws.Range("L86").Value = 5899.9
$ws.Range("M86").Value = -142865867
$ws.Range("N86").Value = -8145.9

$ws.Range("H89").Value = 58831052
$ws.Range("I89").Value = 142866990
$ws.Range("J89").Value = 5899.9
$ws.Range("K89").Value = 714334950
$ws.Range("L89").Value = 29499.5
$ws.Range("M89").Value = -714329334
$ws.Range("N89").Value = -40731.5

$ws.Range("H105").Value = 2997.5
$ws.Range("I105").Value = 2995
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2995
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -1248
$ws.Range("N105").Value = -6494

$ws.Range("H132").Value = 3799.5715
$ws.Range("I132").Value = 2202.6667
$ws.Range("J132").Value = 4997.25
$ws.Range("K132").Value = 6608.000100000001
$ws.Range("L132").Value = 14991.75
$ws.Range("M132").Value = -4078.000100000001
$ws.Range("N132").Value = -20051.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 11360.556
$ws.Range("I23").Value = 94.666664
$ws.Range("J23").Value = 16993.5
$ws.Range("K23").Value = 283.999992
$ws.Range("L23").Value = 50980.5
$ws.Range("M23").Value = -48.99999200000002
$ws.Range("N23").Value = -51450.5

$ws.Range("H113").Value = 579.1277
$ws.Range("I113").Value = 490.58823
$ws.Range("J113").Value = 629.3
$ws.Range("K113").Value = 1471.76469
$ws.Range("L113").Value = 1887.9
$ws.Range("M113").Value = 698.23531
$ws.Range("N113").Value = -6227.9

$ws.Range("H114").Value = 526.5
$ws.Range("I114").Value = 198.0625
$ws.Range("J114").Value = 1840.25
$ws.Range("K114").Value = 594.1875
$ws.Range("L114").Value = 5520.75
$ws.Range("M114").Value = 2659.8125
$ws.Range("N114").Value = -12028.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2265.5557
$ws.Range("I97").Value = 2298.75
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 2298.75
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -1802.75
$ws.Range("N97").Value = -2992

$ws.Range("H126").Value = 10419952
$ws.Range("I126").Value = 16670257
$ws.Range("J126").Value = 2779.111
$ws.Range("K126").Value = 50010771
$ws.Range("L126").Value = 8337.332999999999
$ws.Range("M126").Value = -50008301
$ws.Range("N126").Value = -13277.333

$ws.Range("H132").Value = 5255.5
$ws.Range("I132").Value = 5512
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 16536
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -14006
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 58827624
$ws.Range("I40").Value = 111114120
$ws.Range("J40").Value = 5311.875
$ws.Range("K40").Value = 111114120
$ws.Range("L40").Value = 5311.875
$ws.Range("M40").Value = -111113984
$ws.Range("N40").Value = -5583.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6057.75
$ws.Range("I96").Value = 2509.3333
$ws.Range("J96").Value = 10620
$ws.Range("K96").Value = 2509.3333
$ws.Range("L96").Value = 10620
$ws.Range("M96").Value = -1136.3333
$ws.Range("N96").Value = -13366

$ws.Range("H132").Value = 5982.2856
$ws.Range("I132").Value = 5626
$ws.Range("K132").Value = 16878
$ws.Range("M132").Value = -14348
